$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.441.61"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "3.347.34"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.70"
$ws.Range("E5").Value = "  +4.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "559.78"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.340.19"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.182"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.585"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.92"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "3.882.26"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "599.65"
$ws.Range("E16").Value = "  -6.22%  "
$ws.Range("D17").Value = "66.484.51"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.03"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.118"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.336.56"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.903"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.62"
$ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.96"
$ws.Range("E25").Value = "  -7.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.99"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.76"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.70"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.70"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.72"
$ws.Range("E32").Value = "  +5.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.95"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "581.06"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.03"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "3.706.52"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.48"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.88"
$ws.Range("E41").Value = "  +5.09%  "
$ws.Range("D42").Value = "0.0₃0710"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.23"
$ws.Range("E43").Value = "  -8.16%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.127"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").Value = "  +5.59%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.343"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0421"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.58"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("E51").Value = "  -0.09%  "
